$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 36 (they will become rows 37 and 38).
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(38).Insert()

# Capture the OLD (pre-edit) values of rows 35 and 36 before they get overwritten.
$oldRow35 = $ws.Range("A35:R35").Value2
$oldRow36 = $ws.Range("A36:R36").Value2

# Copy the OLD values into the new rows 37 and 38.
$ws.Range("A37:R37").Value2 = $oldRow35
$ws.Range("A38:R38").Value2 = $oldRow36

# Carry over the date number format to the new rows' D column.
$ws.Range("D37").NumberFormat = $ws.Range("D35").NumberFormat
$ws.Range("D38").NumberFormat = $ws.Range("D35").NumberFormat

# Now update row 35 with the new week's values.
$ws.Range("D35").Value2 = 45142
$ws.Range("H35").Value2 = "Española"
$ws.Range("J35").Value2 = 560
$ws.Range("N35").Value2 = "$/caja 30 unidades"
$ws.Range("O35").Value2 = "Provincia del Elquí"
$ws.Range("P35").Value2 = 483
$ws.Range("Q35").Value2 = 30

# Now update row 36 with the new week's values.
$ws.Range("D36").Value2 = 45142
$ws.Range("J36").Value2 = 600
$ws.Range("K36").Value2 = 9000
$ws.Range("L36").Value2 = 10000
$ws.Range("M36").Value2 = 9500
$ws.Range("P36").Value2 = 238
